$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'sorgen',
    'liegen',
    'filmen',
    'sichern',
    'drücken',
    'rufen',
    'enden',
    'freuen',
    'klettern',
    'arten',
    'fallen',
    'pressen',
    'spenden',
    'runden',
    'rücken',
    'heilen',
    'zögern',
    'fehlen',
    'zeigen',
    'spinnen',
    'schenken',
    'tollen',
    'biegen',
    'melden',
    'fügen',
    'pflegen',
    'schwächen',
    'ehren',
    'kosten',
    'reizen',
    'husten',
    'rasen',
    'heben',
    'trauen',
    'platzen',
    'lächeln',
    'malen',
    'reisen',
    'werden',
    'backen',
    'trennen',
    'landen',
    'werfen',
    'hören',
    'kämpfen',
    'sperren',
    'spannen',
    'betteln',
    'trotzen',
    'stimmen',
    'töten',
    'klingen',
    'wecken',
    'lenken',
    'wellen',
    'planen',
    'sprengen',
    'scheinen',
    'schulden',
    'duschen',
    'leihen',
    'schützen',
    'dringen',
    'jubeln',
    'fließen',
    'wagen',
    'sterben',
    'atmen',
    'tropfen',
    'helfen',
    'spüren',
    'beißen',
    'mauern',
    'passen',
    'wollen',
    'parken',
    'pflanzen',
    'saufen',
    'formen',
    'zielen',
    'wirken',
    'siegen',
    'sagen',
    'kümmern',
    'klagen',
    'schwingen',
    'treiben',
    'schwören',
    'schmecken',
    'bluten',
    'schlafen',
    'kennen',
    'schrecken',
    'teilen',
    'orten',
    'lockern',
    'spielen',
    'gelten',
    'feiern',
    'gleichen',
    'irren',
    'bitten',
    'knarren',
    'proben',
    'zahlen',
    'folgen',
    'liefern',
    'graben',
    'meistern',
    'stehlen',
    'bergen',
    'warnen',
    'heulen',
    'führen',
    'treten',
    'ändern',
    'boxen',
    'zünden',
    'achten',
    'segnen',
    'dienen',
    'sinken',
    'lassen',
    'zeichnen',
    'weichen',
    'äußern',
    'eignen',
    'heißen',
    'pfeifen',
    'ärgern',
    'rechnen',
    'streichen',
    'wehtun',
    'flüchten',
    'stechen',
    'erben',
    'wundern',
    'rauchen',
    'decken',
    'wüten',
    'scheitern',
    'mühen',
    'stören',
    'dürfen',
    'suchen',
    'machen',
    'wenden',
    'buchen',
    'münzen',
    'brauchen',
    'fällen',
    'bauen',
    'lügen',
    'seufzen',
    'handeln',
    'drehen',
    'fangen',
    'stammen',
    'nutzen',
    'hindern',
    'jagen',
    'geben',
    'altern',
    'wechseln',
    'wachsen',
    'greifen',
    'bellen',
    'fischen',
    'schalten',
    'feuern',
    'stecken',
    'narren',
    'fahren',
    'wehren',
    'schreiten',
    'prügeln',
    'schlucken',
    'räumen',
    'loben',
    'flehen',
    'mögen',
    'grüßen',
    'hacken',
    'kehren',
    'lesen',
    'kichern',
    'morden',
    'hauen',
    'gründen',
    'drohen',
    'doppeln',
    'quälen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}
